$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2")

# Update Michael's birth month from April to June
$ws.Range("B2").Value = "June"

# Remove Daniel's row (originally row 6) and John's row (originally row 3).
# Delete the lower row first so the earlier row index stays valid.
$ws.Rows("6").Delete() | Out-Null
$ws.Rows("3").Delete() | Out-Null

# Update the active selection to match the saved view
$ws.Range("E14").Select() | Out-Null
